$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 45: "10. Regular Expression Matching" (finished, not yet reviewed)
$ws.Rows.Item(45).RowHeight = 42

# --- A45: title ---
$a = $ws.Range("A45")
$a.Value = "10. Regular Expression Matching"
$a.Font.Name = "Times New Roman"
$a.Font.Size = 11
$a.Interior.Color = 255
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4108
$a.WrapText = $true

# --- B45: difficulty ---
$b = $ws.Range("B45")
$b.Value = "Hard"
$b.Font.Name = "Times New Roman"
$b.Font.Size = 11
$b.Interior.Color = 255
$b.HorizontalAlignment = -4108
$b.VerticalAlignment = -4108

# --- C45: link (hyperlink) ---
$c = $ws.Range("C45")
$c.Value = "https://leetcode.com/problems/regular-expression-matching/"
$ws.Hyperlinks.Add($c, "https://leetcode.com/problems/regular-expression-matching/") | Out-Null
$c.Interior.Color = 255
$c.WrapText = $true

# --- D45: date finished ---
$d = $ws.Range("D45")
$d.Value = 44548
$d.NumberFormat = "mm-dd-yy"
$d.Font.Name = "Times New Roman"
$d.Font.Size = 11
$d.Interior.Color = 255
$d.HorizontalAlignment = -4108
$d.VerticalAlignment = -4108

# --- E45: category ---
$e = $ws.Range("E45")
$e.Value = "动态规划，字符串编辑"
$e.Font.Name = "宋体"
$e.Font.Size = 11
$e.Font.Family = 3
$e.Interior.Color = 255
$e.HorizontalAlignment = -4108
$e.VerticalAlignment = -4108

# --- F45: notes (rich text, 3 runs) ---
$f = $ws.Range("F45")
$full = "分类讨论，关键就是带星号的情况下，是否要重复，状态转移方程不同；注意等价dp[0][0]情况的初始化"
$f.Value = $full
$f.Font.Name = "Times New Roman"
$f.Font.Size = 11
$f.Font.Family = 3
$f.Interior.Color = 255
$f.HorizontalAlignment = -4108
$f.VerticalAlignment = -4108
$f.WrapText = $true

$f.Characters(1, 36).Font.Name = "宋体"
$f.Characters(1, 36).Font.Size = 11
$f.Characters(1, 36).Font.Family = 3

$f.Characters(37, 8).Font.Name = "Times New Roman"
$f.Characters(37, 8).Font.Size = 11
$f.Characters(37, 8).Font.Family = 1

$f.Characters(45, 6).Font.Name = "宋体"
$f.Characters(45, 6).Font.Size = 11
$f.Characters(45, 6).Font.Family = 3

# --- G45: review state (not yet reviewed) ---
$g = $ws.Range("G45")
$g.Value = "未复习"
$g.Font.Name = "宋体"
$g.Font.Size = 11
$g.Font.Family = 3
$g.HorizontalAlignment = -4108
$g.VerticalAlignment = -4108

# --- H45: mastered marker ---
$h = $ws.Range("H45")
$h.Value = "⭕"
$h.Font.Name = "宋体"
$h.Font.Size = 11
$h.Font.Bold = $true
$h.Font.Family = 3
$h.HorizontalAlignment = -4108
$h.VerticalAlignment = -4108

# Keep the view roughly where the author left it
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("H46").Select()
